$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-03-29 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-03-30 Thursday", 2) | Out-Null
$d.Content.Find.Execute("69×60=", $true, $false, $false, $false, $false, $true, 1, $false, "79×51=", 2) | Out-Null
$d.Content.Find.Execute("95×33=", $true, $false, $false, $false, $false, $true, 1, $false, "44×67=", 2) | Out-Null
$d.Content.Find.Execute("94×12=", $true, $false, $false, $false, $false, $true, 1, $false, "93×21=", 2) | Out-Null
$d.Content.Find.Execute("78×41=", $true, $false, $false, $false, $false, $true, 1, $false, "10×13=", 2) | Out-Null
$d.Content.Find.Execute("58×34=", $true, $false, $false, $false, $false, $true, 1, $false, "58×97=", 2) | Out-Null
$d.Content.Find.Execute("12×49=", $true, $false, $false, $false, $false, $true, 1, $false, "78×93=", 2) | Out-Null
$d.Content.Find.Execute("21×15=", $true, $false, $false, $false, $false, $true, 1, $false, "94×34=", 2) | Out-Null
$d.Content.Find.Execute("79×55=", $true, $false, $false, $false, $false, $true, 1, $false, "97×31=", 2) | Out-Null
$d.Content.Find.Execute("43×20=", $true, $false, $false, $false, $false, $true, 1, $false, "29×82=", 2) | Out-Null
$d.Content.Find.Execute("22×74=", $true, $false, $false, $false, $false, $true, 1, $false, "38×46=", 2) | Out-Null
$d.Content.Find.Execute("90×24=", $true, $false, $false, $false, $false, $true, 1, $false, "17×92=", 2) | Out-Null
$d.Content.Find.Execute("16×42=", $true, $false, $false, $false, $false, $true, 1, $false, "68×27=", 2) | Out-Null
$d.Content.Find.Execute("18×50=", $true, $false, $false, $false, $false, $true, 1, $false, "68×34=", 2) | Out-Null
$d.Content.Find.Execute("95×76=", $true, $false, $false, $false, $false, $true, 1, $false, "95×54=", 2) | Out-Null
$d.Content.Find.Execute("89×53=", $true, $false, $false, $false, $false, $true, 1, $false, "13×33=", 2) | Out-Null
$d.Content.Find.Execute("79×49=", $true, $false, $false, $false, $false, $true, 1, $false, "34×36=", 2) | Out-Null
$d.Content.Find.Execute("48×32=", $true, $false, $false, $false, $false, $true, 1, $false, "67×88=", 2) | Out-Null
$d.Content.Find.Execute("11×60=", $true, $false, $false, $false, $false, $true, 1, $false, "12×19=", 2) | Out-Null
$d.Content.Find.Execute("89×11=", $true, $false, $false, $false, $false, $true, 1, $false, "11×79=", 2) | Out-Null
$d.Content.Find.Execute("40×85=", $true, $false, $false, $false, $false, $true, 1, $false, "64×69=", 2) | Out-Null
$d.Content.Find.Execute("30×11=", $true, $false, $false, $false, $false, $true, 1, $false, "39×44=", 2) | Out-Null
$d.Content.Find.Execute("71×66=", $true, $false, $false, $false, $false, $true, 1, $false, "74×58=", 2) | Out-Null
$d.Content.Find.Execute("83×39=", $true, $false, $false, $false, $false, $true, 1, $false, "87×55=", 2) | Out-Null
$d.Content.Find.Execute("51×20=", $true, $false, $false, $false, $false, $true, 1, $false, "12×59=", 2) | Out-Null
$d.Content.Find.Execute("84×91=", $true, $false, $false, $false, $false, $true, 1, $false, "81×84=", 2) | Out-Null
$d.Content.Find.Execute("85×92=", $true, $false, $false, $false, $false, $true, 1, $false, "91×36=", 2) | Out-Null
$d.Content.Find.Execute("10×34=", $true, $false, $false, $false, $false, $true, 1, $false, "98×64=", 2) | Out-Null
$d.Content.Find.Execute("76×86=", $true, $false, $false, $false, $false, $true, 1, $false, "99×71=", 2) | Out-Null
$d.Content.Find.Execute("76×40=", $true, $false, $false, $false, $false, $true, 1, $false, "50×18=", 2) | Out-Null
$d.Content.Find.Execute("34×53=", $true, $false, $false, $false, $false, $true, 1, $false, "87×30=", 2) | Out-Null
$d.Content.Find.Execute("33×36=", $true, $false, $false, $false, $false, $true, 1, $false, "26×77=", 2) | Out-Null
$d.Content.Find.Execute("17×52=", $true, $false, $false, $false, $false, $true, 1, $false, "96×75=", 2) | Out-Null
$d.Content.Find.Execute("94×67=", $true, $false, $false, $false, $false, $true, 1, $false, "65×82=", 2) | Out-Null
$d.Content.Find.Execute("87×68=", $true, $false, $false, $false, $false, $true, 1, $false, "64×86=", 2) | Out-Null
$d.Content.Find.Execute("32×68=", $true, $false, $false, $false, $false, $true, 1, $false, "99×51=", 2) | Out-Null
$d.Content.Find.Execute("85×21=", $true, $false, $false, $false, $false, $true, 1, $false, "33×94=", 2) | Out-Null
$d.Content.Find.Execute("87×94=", $true, $false, $false, $false, $false, $true, 1, $false, "53×48=", 2) | Out-Null
$d.Content.Find.Execute("100×12=", $true, $false, $false, $false, $false, $true, 1, $false, "100×65=", 2) | Out-Null
$d.Content.Find.Execute("48×68=", $true, $false, $false, $false, $false, $true, 1, $false, "70×73=", 2) | Out-Null
$d.Content.Find.Execute("97×57=", $true, $false, $false, $false, $false, $true, 1, $false, "81×68=", 2) | Out-Null
$d.Content.Find.Execute("35×67=", $true, $false, $false, $false, $false, $true, 1, $false, "98×49=", 2) | Out-Null
$d.Content.Find.Execute("16×46=", $true, $false, $false, $false, $false, $true, 1, $false, "74×77=", 2) | Out-Null
$d.Content.Find.Execute("29×99=", $true, $false, $false, $false, $false, $true, 1, $false, "63×69=", 2) | Out-Null
$d.Content.Find.Execute("75×92=", $true, $false, $false, $false, $false, $true, 1, $false, "61×92=", 2) | Out-Null
$d.Content.Find.Execute("34×26=", $true, $false, $false, $false, $false, $true, 1, $false, "86×95=", 2) | Out-Null
$d.Content.Find.Execute("79×87=", $true, $false, $false, $false, $false, $true, 1, $false, "73×65=", 2) | Out-Null
$d.Content.Find.Execute("36×78=", $true, $false, $false, $false, $false, $true, 1, $false, "34×29=", 2) | Out-Null
$d.Content.Find.Execute("84×40=", $true, $false, $false, $false, $false, $true, 1, $false, "72×21=", 2) | Out-Null
$d.Content.Find.Execute("14×99=", $true, $false, $false, $false, $false, $true, 1, $false, "64×48=", 2) | Out-Null
$d.Content.Find.Execute("92×51=", $true, $false, $false, $false, $false, $true, 1, $false, "39×89=", 2) | Out-Null
$d.Content.Find.Execute("33×11=", $true, $false, $false, $false, $false, $true, 1, $false, "44×70=", 2) | Out-Null
$d.Content.Find.Execute("47×42=", $true, $false, $false, $false, $false, $true, 1, $false, "55×29=", 2) | Out-Null
$d.Content.Find.Execute("66×46=", $true, $false, $false, $false, $false, $true, 1, $false, "52×44=", 2) | Out-Null
$d.Content.Find.Execute("78×43=", $true, $false, $false, $false, $false, $true, 1, $false, "99×15=", 2) | Out-Null
$d.Content.Find.Execute("17×29=", $true, $false, $false, $false, $false, $true, 1, $false, "58×25=", 2) | Out-Null
$d.Content.Find.Execute("40×91=", $true, $false, $false, $false, $false, $true, 1, $false, "96×62=", 2) | Out-Null
$d.Content.Find.Execute("59×100=", $true, $false, $false, $false, $false, $true, 1, $false, "83×16=", 2) | Out-Null
$d.Content.Find.Execute("83×84=", $true, $false, $false, $false, $false, $true, 1, $false, "35×38=", 2) | Out-Null
$d.Content.Find.Execute("63×50=", $true, $false, $false, $false, $false, $true, 1, $false, "67×21=", 2) | Out-Null
$d.Content.Find.Execute("48×72=", $true, $false, $false, $false, $false, $true, 1, $false, "54×57=", 2) | Out-Null
$d.Content.Find.Execute("60×97=", $true, $false, $false, $false, $false, $true, 1, $false, "78×31=", 2) | Out-Null
$d.Content.Find.Execute("42×97=", $true, $false, $false, $false, $false, $true, 1, $false, "57×72=", 2) | Out-Null
$d.Content.Find.Execute("31×19=", $true, $false, $false, $false, $false, $true, 1, $false, "80×12=", 2) | Out-Null
$d.Content.Find.Execute("57×88=", $true, $false, $false, $false, $false, $true, 1, $false, "56×13=", 2) | Out-Null
$d.Content.Find.Execute("73×34=", $true, $false, $false, $false, $false, $true, 1, $false, "58×43=", 2) | Out-Null
$d.Content.Find.Execute("83×62=", $true, $false, $false, $false, $false, $true, 1, $false, "72×22=", 2) | Out-Null
$d.Content.Find.Execute("89×42=", $true, $false, $false, $false, $false, $true, 1, $false, "13×48=", 2) | Out-Null
$d.Content.Find.Execute("25×20=", $true, $false, $false, $false, $false, $true, 1, $false, "65×90=", 2) | Out-Null
$d.Content.Find.Execute("25×53=", $true, $false, $false, $false, $false, $true, 1, $false, "92×70=", 2) | Out-Null
$d.Content.Find.Execute("42×14=", $true, $false, $false, $false, $false, $true, 1, $false, "100×71=", 2) | Out-Null
$d.Content.Find.Execute("22×48=", $true, $false, $false, $false, $false, $true, 1, $false, "45×49=", 2) | Out-Null
$d.Content.Find.Execute("40×50=", $true, $false, $false, $false, $false, $true, 1, $false, "51×39=", 2) | Out-Null
$d.Content.Find.Execute("100×45=", $true, $false, $false, $false, $false, $true, 1, $false, "49×85=", 2) | Out-Null
$d.Content.Find.Execute("80×71=", $true, $false, $false, $false, $false, $true, 1, $false, "28×63=", 2) | Out-Null
$d.Content.Find.Execute("70×30=", $true, $false, $false, $false, $false, $true, 1, $false, "76×79=", 2) | Out-Null
$d.Content.Find.Execute("14×50=", $true, $false, $false, $false, $false, $true, 1, $false, "89×41=", 2) | Out-Null
$d.Content.Find.Execute("53×75=", $true, $false, $false, $false, $false, $true, 1, $false, "39×75=", 2) | Out-Null
$d.Content.Find.Execute("83×55=", $true, $false, $false, $false, $false, $true, 1, $false, "46×50=", 2) | Out-Null
$d.Content.Find.Execute("85×80=", $true, $false, $false, $false, $false, $true, 1, $false, "78×63=", 2) | Out-Null
$d.Content.Find.Execute("19×43=", $true, $false, $false, $false, $false, $true, 1, $false, "71×98=", 2) | Out-Null
$d.Content.Find.Execute("57×48=", $true, $false, $false, $false, $false, $true, 1, $false, "14×76=", 2) | Out-Null
$d.Content.Find.Execute("69×48=", $true, $false, $false, $false, $false, $true, 1, $false, "16×38=", 2) | Out-Null
$d.Content.Find.Execute("51×77=", $true, $false, $false, $false, $false, $true, 1, $false, "83×38=", 2) | Out-Null
$d.Content.Find.Execute("65×88=", $true, $false, $false, $false, $false, $true, 1, $false, "91×17=", 2) | Out-Null
$d.Content.Find.Execute("100×48=", $true, $false, $false, $false, $false, $true, 1, $false, "82×39=", 2) | Out-Null
$d.Content.Find.Execute("18×82=", $true, $false, $false, $false, $false, $true, 1, $false, "99×51=", 2) | Out-Null
$d.Content.Find.Execute("31×90=", $true, $false, $false, $false, $false, $true, 1, $false, "55×86=", 2) | Out-Null
$d.Content.Find.Execute("27×18=", $true, $false, $false, $false, $false, $true, 1, $false, "53×81=", 2) | Out-Null
$d.Content.Find.Execute("22×88=", $true, $false, $false, $false, $false, $true, 1, $false, "33×31=", 2) | Out-Null
$d.Content.Find.Execute("33×68=", $true, $false, $false, $false, $false, $true, 1, $false, "23×41=", 2) | Out-Null
$d.Content.Find.Execute("61×95=", $true, $false, $false, $false, $false, $true, 1, $false, "32×30=", 2) | Out-Null
$d.Content.Find.Execute("66×62=", $true, $false, $false, $false, $false, $true, 1, $false, "29×45=", 2) | Out-Null
$d.Content.Find.Execute("91×20=", $true, $false, $false, $false, $false, $true, 1, $false, "15×70=", 2) | Out-Null
$d.Content.Find.Execute("55×24=", $true, $false, $false, $false, $false, $true, 1, $false, "15×56=", 2) | Out-Null
$d.Content.Find.Execute("51×24=", $true, $false, $false, $false, $false, $true, 1, $false, "47×87=", 2) | Out-Null
$d.Content.Find.Execute("63×95=", $true, $false, $false, $false, $false, $true, 1, $false, "91×32=", 2) | Out-Null
$d.Content.Find.Execute("15×63=", $true, $false, $false, $false, $false, $true, 1, $false, "72×83=", 2) | Out-Null
$d.Content.Find.Execute("75×61=", $true, $false, $false, $false, $false, $true, 1, $false, "52×10=", 2) | Out-Null
$d.Content.Find.Execute("91×58=", $true, $false, $false, $false, $false, $true, 1, $false, "16×68=", 2) | Out-Null
$d.Content.Find.Execute("90×30=", $true, $false, $false, $false, $false, $true, 1, $false, "61×12=", 2) | Out-Null
